# =====================================================================
# Edit script: "break out stock.yaml completed"
#
# 1. Resets a batch of existing weekly rows "Q"/"R" helper-flag columns
#    (detect_structure / backup) from stale 1/2 markers back to 0, and
#    bumps O1134 (isPivot) from 0 to 3.
# 2. Normalises the previously-blank R1136/R1137 "backup" cells to an
#    explicit numeric 0 (they used to be empty placeholders).
# 3. Appends 22 new weekly OHLCV rows (1138-1159), continuing the
#    Datetime series from 2024-07-01 through 2024-11-25, with their
#    "backup" column left blank exactly like the most-recent existing
#    rows were before being finalised.
# =====================================================================

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) & 2): point updates on existing cells -------------------------
$cellUpdates = @(
    @("Q56", 0),
    @("R56", 0),
    @("Q59", 0),
    @("R59", 0),
    @("Q64", 0),
    @("R64", 0),
    @("Q68", 0),
    @("R68", 0),
    @("R76", 0),
    @("R79", 0),
    @("R92", 0),
    @("R101", 0),
    @("R108", 0),
    @("R114", 0),
    @("R125", 0),
    @("R128", 0),
    @("R136", 0),
    @("R163", 0),
    @("R178", 0),
    @("R180", 0),
    @("R197", 0),
    @("R201", 0),
    @("R207", 0),
    @("R217", 0),
    @("R221", 0),
    @("R228", 0),
    @("R240", 0),
    @("R264", 0),
    @("R273", 0),
    @("R283", 0),
    @("R286", 0),
    @("R294", 0),
    @("R301", 0),
    @("R311", 0),
    @("R324", 0),
    @("R326", 0),
    @("R334", 0),
    @("R345", 0),
    @("R364", 0),
    @("R368", 0),
    @("R382", 0),
    @("R385", 0),
    @("R394", 0),
    @("R407", 0),
    @("R409", 0),
    @("R430", 0),
    @("R436", 0),
    @("R446", 0),
    @("R447", 0),
    @("R453", 0),
    @("R457", 0),
    @("R478", 0),
    @("R491", 0),
    @("R499", 0),
    @("R514", 0),
    @("R519", 0),
    @("R526", 0),
    @("R544", 0),
    @("R558", 0),
    @("R562", 0),
    @("R578", 0),
    @("R581", 0),
    @("R594", 0),
    @("R602", 0),
    @("R614", 0),
    @("R619", 0),
    @("R624", 0),
    @("R639", 0),
    @("R651", 0),
    @("R666", 0),
    @("R681", 0),
    @("R688", 0),
    @("R696", 0),
    @("R709", 0),
    @("R714", 0),
    @("R729", 0),
    @("R740", 0),
    @("R745", 0),
    @("R747", 0),
    @("R752", 0),
    @("R758", 0),
    @("R770", 0),
    @("R778", 0),
    @("R782", 0),
    @("R784", 0),
    @("R795", 0),
    @("R819", 0),
    @("R824", 0),
    @("R828", 0),
    @("R833", 0),
    @("R840", 0),
    @("R847", 0),
    @("R850", 0),
    @("R854", 0),
    @("R859", 0),
    @("R864", 0),
    @("R870", 0),
    @("R883", 0),
    @("R904", 0),
    @("R919", 0),
    @("R929", 0),
    @("R934", 0),
    @("R944", 0),
    @("R950", 0),
    @("R958", 0),
    @("R968", 0),
    @("R977", 0),
    @("R980", 0),
    @("R986", 0),
    @("R995", 0),
    @("R1006", 0),
    @("R1019", 0),
    @("R1027", 0),
    @("R1032", 0),
    @("R1039", 0),
    @("R1063", 0),
    @("R1078", 0),
    @("R1083", 0),
    @("R1091", 0),
    @("R1100", 0),
    @("R1105", 0),
    @("R1108", 0),
    @("R1131", 0),
    @("O1134", 3)
)
foreach ($u in $cellUpdates) {
    $ws.Range($u[0]).Value = $u[1]
}

# --- 3) append new rows 1138-1159 --------------------------------------
# Each tuple: (row, A:Datetime .. Q:detect_structure, R:backup-or-null)
$newRows = @(
    @(1138, 45474, 137.3899993896484, 137.3999938964844, 133.5, 135.7299957275391, 132.2348022460938, 60453597, 2024, 7, 1, 0, 0, 0, 27, 0, 0, 0, $null),
    @(1139, 45481, 135, 141.0700073242188, 133, 136.1100006103516, 132.6050262451172, 63335796, 2024, 7, 8, 0, 0, 0, 28, 0, 0, 0, $null),
    @(1140, 45488, 136.5, 142.7400054931641, 135.1000061035156, 135.6499938964844, 132.1568603515625, 59386920, 2024, 7, 15, 0, 0, 0, 29, 0, 0, 0, $null),
    @(1141, 45495, 134.6499938964844, 138.9600067138672, 129.1000061035156, 132.8699951171875, 132.8699951171875, 51816593, 2024, 7, 22, 0, 0, 0, 30, 0, 0, 0, $null),
    @(1142, 45502, 134.25, 137.4600067138672, 131.7899932861328, 133.2700042724609, 133.2700042724609, 50631659, 2024, 7, 29, 0, 0, 0, 31, 0, 0, 0, $null),
    @(1143, 45509, 129.5, 129.8399963378906, 120.0999984741211, 123.0199966430664, 123.0199966430664, 74062845, 2024, 8, 5, 0, 0, 0, 32, 0, 0, 0, $null),
    @(1144, 45516, 122.5100021362305, 122.879997253418, 116.3000030517578, 117.5299987792969, 117.5299987792969, 47374516, 2024, 8, 12, 0, 0, 0, 33, 2, 0, 0, $null),
    @(1145, 45523, 118.5100021362305, 128.0200042724609, 118.2099990844727, 127.0800018310547, 127.0800018310547, 51119901, 2024, 8, 19, 0, 0, 0, 34, 0, 0, 0, $null),
    @(1146, 45530, 127.0999984741211, 127.75, 121.0999984741211, 121.5400009155273, 121.5400009155273, 44504757, 2024, 8, 26, 0, 0, 0, 35, 0, 0, 0, $null),
    @(1147, 45537, 122.4899978637695, 123.620002746582, 120, 121.1999969482422, 121.1999969482422, 38306091, 2024, 9, 2, 0, 0, 0, 36, 0, 0, 0, $null),
    @(1148, 45544, 120.5100021362305, 123, 117.1100006103516, 120.1900024414062, 120.1900024414062, 45179100, 2024, 9, 9, 0, 0, 0, 37, 0, 0, 0, $null),
    @(1149, 45551, 122.8000030517578, 127.3000030517578, 119.3000030517578, 123.4300003051758, 123.4300003051758, 66453242, 2024, 9, 16, 0, 0, 0, 38, 0, 0, 0, $null),
    @(1150, 45558, 124.4000015258789, 129, 121.4700012207031, 123.5100021362305, 123.5100021362305, 75160673, 2024, 9, 23, 0, 0, 0, 39, 1, 0, 1, $null),
    @(1151, 45565, 123.1999969482422, 123.9000015258789, 116.7300033569336, 118.8600006103516, 118.8600006103516, 41807884, 2024, 9, 30, 0, 0, 0, 40, 0, 0, 0, $null),
    @(1152, 45572, 119.4499969482422, 120.1900024414062, 112.5899963378906, 114.120002746582, 114.120002746582, 37664068, 2024, 10, 7, 0, 0, 0, 41, 0, 0, 0, $null),
    @(1153, 45579, 114.4499969482422, 116.4499969482422, 110.2699966430664, 112.2200012207031, 112.2200012207031, 50882606, 2024, 10, 14, 0, 0, 0, 42, 0, 0, 0, $null),
    @(1154, 45586, 112.6500015258789, 114.1100006103516, 106.6800003051758, 108.2399978637695, 108.2399978637695, 55752622, 2024, 10, 21, 0, 0, 0, 43, 2, 0, 0, $null),
    @(1155, 45593, 108.7900009155273, 118.6999969482422, 108.1600036621094, 117.8499984741211, 117.8499984741211, 36244424, 2024, 10, 28, 0, 0, 0, 44, 0, 0, 0, $null),
    @(1156, 45600, 118.6999969482422, 122.4100036621094, 113.1600036621094, 117.5199966430664, 117.5199966430664, 42323420, 2024, 11, 4, 0, 0, 0, 45, 0, 0, 2, $null),
    @(1157, 45607, 117.5, 120.6399993896484, 113.25, 113.9000015258789, 113.9000015258789, 29437081, 2024, 11, 11, 0, 0, 0, 46, 0, 0, 0, $null),
    @(1158, 45614, 114, 117.7900009155273, 112, 115.0400009155273, 115.0400009155273, 29074206, 2024, 11, 18, 0, 0, 0, 47, 0, 0, 0, $null),
    @(1159, 45621, 117.0899963378906, 126.1999969482422, 117.0899963378906, 121.620002746582, 121.620002746582, 67864331, 2024, 11, 25, 0, 0, 0, 48, 0, 0, 0, $null)
)

foreach ($row in $newRows) {
    $rownum = $row[0]
    # Column A: Datetime -- carries the same date/time style as the rest of the column
    $ws.Cells.Item($rownum, 1).Value = $row[1]
    $ws.Cells.Item($rownum, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
    # Columns B (Open) .. Q (detect_structure)
    for ($c = 2; $c -le 17; $c++) {
        $ws.Cells.Item($rownum, $c).Value = $row[$c]
    }
    # Column R (backup) stays blank for brand-new rows, same as the tail
    # of the sheet looked like before being backfilled.
    if ($row[18] -ne $null) {
        $ws.Cells.Item($rownum, 18).Value2 = $row[18]
    } else {
        $ws.Cells.Item($rownum, 18).Value2 = ""
    }
}
